# Apply edits described by the diff:
# - F1 changes from "Caretaker" to "User"
# - New row 3 added: Juan | Carlos | banan | 4567 | Underage_user | User
# - Selection moves to E5
# - Column E gets a custom width

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update F1 value
$ws.Range("F1").Value = "User"

# Add new row 3
$ws.Range("A3").Value = "Juan"
$ws.Range("B3").Value = "Carlos"
$ws.Range("C3").Value = "banan"
$ws.Range("D3").Value = 4567
$ws.Range("E3").Value = "Underage_user"
$ws.Range("F3").Value = "User"

# Set column E width to match target (stored XML width 17.54296875 "characters").
# The engine quantizes ColumnWidth to whole pixels (1/6-character steps at the
# default 6px digit width), so the nearest reachable value is used here.
$ws.Columns.Item(5).ColumnWidth = 16.6666666666667

# Update selection to E5
$ws.Range("E5").Select()
